# Word COM-interop edit script.
#
# Target change (per the supplied OOXML diff):
#   1. Insert a new, completely empty paragraph ("<w:p/>") right after the
#      "A fourth paragraph?????" paragraph, i.e. immediately before the
#      "Yea yea, this is an conclusion" paragraph — the author was trying
#      to give the two "Yea yea..." blocks (intro/conclusion) a visual
#      separator.
#   2. In that "conclusion" paragraph, the grammar-check run break around
#      the word "an" (<w:proofErr w:type="gramStart"/> ... "an" ...
#      <w:proofErr w:type="gramEnd"/>) gets collapsed: the three runs
#      ", this is " + "an" + " conclusion" become one single run
#      ", this is an conclusion" and the gramStart/gramEnd proofErr
#      markers disappear.

$d = $word.ActiveDocument

# --- 1. Insert a blank paragraph after "A fourth paragraph?????" --------
# Using Find/Replace with the "^p" paragraph-mark wildcard (rather than
# Range.InsertParagraphAfter) is what produces a bare "<w:p/>" element
# instead of a paragraph carrying a leftover empty run.
$d.Content.Find.Execute(
    "A fourth paragraph?????", $true, $false, $false, $false, $false,
    $true, 1, $false, "A fourth paragraph?????^p", 2) | Out-Null

# --- 2. Merge the split "an" runs in the conclusion paragraph -----------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Yea yea, this is an conclusion*") {
        $para.Range.Find.Execute(
            ", this is an conclusion", $true, $false, $false, $false, $false,
            $true, 1, $false, ", this is an conclusion", 2) | Out-Null
        break
    }
}
